$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right count 5 -> 4, Wrong marking -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right total 90 -> 72, Wrong total -5 -> -10, and summary text update
$ws.Range("B12").Value = 72
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "62 / 112"
